$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unmerge existing merged header cells before restructuring rows
$ws.Range("A1:B1").UnMerge()
$ws.Range("A5:B5").UnMerge()
$ws.Range("A9:B9").UnMerge()

# --- Block 1: Assert200_pendingtoapprove (rows 1-4) ---
$ws.Range("A1").Value = "Assert200_pendingtoapprove"
$ws.Range("A2").Value = "GroupStatus"
$ws.Range("B2").Value = "EndPoint"
$ws.Range("A3").Value = "Approved"
$ws.Range("B3").Value = "/conclusionInfo/conclusions"
$ws.Range("A4").Value = $null
$ws.Range("B4").Value = $null

# --- Block 2: Assert200_approved (rows 5-8) ---
$ws.Range("A5").Value = "Assert200_approved"
$ws.Range("A6").Value = "GroupStatus"
$ws.Range("B6").Value = "EndPoint"
$ws.Range("A7").Value = "Approved"
$ws.Range("B7").Value = "/conclusionInfo/conclusions"
$ws.Range("A8").Value = $null
$ws.Range("B8").Value = $null

# --- Block 3: Assert400 (rows 9-11) ---
$ws.Range("A9").Value = "Assert400"
$ws.Range("A10").Value = "GroupStatus"
$ws.Range("B10").Value = "EndPoint"
$ws.Range("A11").Value = "Approved1"
$ws.Range("B11").Value = "/conclusionInfo/conclusions"

# Remove old row 12 content and old empty rows 13-25 that are no longer used
$ws.Rows("12:25").Delete()

# --- Block 4: Assert401 (rows 13-16) ---
$ws.Range("A13").Value = "Assert401"
$ws.Range("A14").Value = "GroupStatus"
$ws.Range("B14").Value = "EndPoint"
$ws.Range("A15").Value = "Approved"
$ws.Range("B15").Value = "/conclusionInfo/conclusions"
$ws.Range("A16").Value = $null

# Apply styles to the new/shifted cells based on source rows (s=6 header, s=5/2 label row, s=3 value row)
$ws.Range("A1:B1").Style = $ws.Range("A9").Style
$ws.Range("A4:B4").Style = $ws.Range("A3").Style
$ws.Range("A5:B5").Style = $ws.Range("A9").Style
$ws.Range("A8:B8").Style = $ws.Range("A7").Style
$ws.Range("A16").Style = $ws.Range("A12").Style

# Re-apply merges
$ws.Range("A1:B1").Merge()
$ws.Range("A5:B5").Merge()
$ws.Range("A9:B9").Merge()
$ws.Range("A13:B13").Merge()

# Selection
$ws.Range("A3").Select()
